$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: add a new sentence after "...proper tables, fields, keys, etc. "
# -----------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("I am also tasked with getting the database setup with the proper tables, fields, keys, etc. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter("With this responsibility, I will be working on making sure the data can be passed back and forth easily.")
}

# -----------------------------------------------------------------
# Change 2: split "...but I may need help with in the next week or two."
# into "...but I may need help with" + " this" + bookmark(_GoBack) + " in the next week or two."
# -----------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("but I may need help with", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter(" this")
}

# Re-find the same anchor text to get a fresh Range and mark the first split point
# (between "...help with" and " this").
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("but I may need help with", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Collapse(0)
    $d.Bookmarks.Add("zzzTempSplit", $rng3)
}

# Mark the second split point (right after the newly inserted " this", before
# " in the next week or two."). This also relocates the single allowed "_GoBack"
# bookmark away from its old location further down in the document.
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("help with this", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $rng4.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng4)
}

# Drop the temporary bookmark now that it has done its job of keeping the
# "with" / " this" runs from re-merging.
if ($d.Bookmarks.Exists("zzzTempSplit")) {
    $d.Bookmarks("zzzTempSplit").Delete()
}

Write-Output "done"
